$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.303.53'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '2.633.67'
$ws.Range("E3").Value = '  -2.66%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'596.61"
$ws.Range("E5").Value = '  -2.06%  '
$ws.Range("D6").Value = "'168.00"
$ws.Range("E6").Value = '  +1.01%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -3.40%  '
$ws.Range("D9").Value = '2.633.00'
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = '  -1.42%  '
$ws.Range("D12").Value = "'0.362"
$ws.Range("E12").Value = '  -0.04%  '
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("E14").Value = '  -1.71%  '
$ws.Range("E15").Value = '  -2.72%  '
$ws.Range("E16").Value = '  -2.65%  '
$ws.Range("D17").Value = '67.252.28'
$ws.Range("E17").Value = '  -1.35%  '
$ws.Range("D18").Value = '2.626.82'
$ws.Range("E18").Value = '  -2.66%  '
$ws.Range("D19").Value = "'12.18"
$ws.Range("E19").Value = '  +3.73%  '
$ws.Range("D20").Value = "'8.11"
$ws.Range("E20").Value = '  +6.34%  '
$ws.Range("D21").Value = "'358.95"
$ws.Range("E21").Value = '  -2.78%  '
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("E23").Value = '  -4.45%  '
$ws.Range("E24").Value = '  +9.23%  '
$ws.Range("E25").Value = '  -5.46%  '
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = "'70.57"
$ws.Range("E27").Value = '  -2.91%  '
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = "'0.0000101"
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("B30").Value = 'Bittensor'
$ws.Range("C30").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D30").Value = "'555.74"
$ws.Range("E30").Value = '  -3.89%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = "'7.94"
$ws.Range("E31").Value = '  -1.99%  '
$ws.Range("B32").Value = 'Fetch.AI'
$ws.Range("C32").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D32").Value = "'1.37"
$ws.Range("E32").Value = '  -2.94%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value = "'1.91"
$ws.Range("E33").Value = '  -3.10%  '
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = "'0.137"
$ws.Range("E34").Value = '  +4.38%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = "'1.51"
$ws.Range("E36").Value = '  -4.97%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = "'157.46"
$ws.Range("E37").Value = '  -0.75%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = "'19.20"
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = "'0.367"
$ws.Range("E39").Value = '  -2.75%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = "'5.19"
$ws.Range("E40").Value = '  -3.51%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = "'1.80"
$ws.Range("E41").Value = '  -2.99%  '
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D42").Value = "'17.94"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = "'2.47"
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").Value = "'40.19"
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = '0.0₆0299'
$ws.Range("E46").Value = '  -2.95%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").Value = "'0.587"
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").Value = "'152.28"
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = "'3.83"
$ws.Range("E49").Value = '  -1.35%  '
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = "'1.73"
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = "'0.0772"
$ws.Range("E51").Value = '  -1.39%  '
